$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2739937603473663
$ws.Range("B1").Value = 0.7964653968811035
$ws.Range("C1").Value = 3.658106565475464
$ws.Range("D1").Value = 3.303115606307983
$ws.Range("E1").Value = 0.8817182183265686
